$wb = $excel.ActiveWorkbook

# --- Hoja1!A1 text update (rates in the daily-conversion note) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.92 = 41634.85 pesos`n✅ 41634.85 pesos = 9.91 = 957.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas sheet numeric updates ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 100.798
$ws2.Range("O10").Value = 4196.71
$ws2.Range("N12").Value = 4200
$ws2.Range("O12").Value = 96.555
